# Apply content edits to sheet "Лист1" of the workbook.
# Only the cells whose text actually changes need to be touched - Excel
# rebuilds/reorders the shared-strings table automatically on save, which
# accounts for the bulk of the index churn visible in the raw XML diff.
# The cells below are written in the same order the corresponding strings
# were newly introduced in the target file, so the regenerated shared
# string table lines up with the target as closely as possible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value  = "By sex"
$ws.Range("C12").Value = "By territory"
$ws.Range("C22").Value = "By age (in month)"
$ws.Range("C28").Value = "Education of mother"
$ws.Range("C29").Value = "Preschool or not /primary"
$ws.Range("C30").Value = "Basic general"
$ws.Range("C31").Value = "Average total"
$ws.Range("C32").Value = "Vocational primary /secondary"
$ws.Range("C33").Value = "Higher"

$ws.Range("C7").Value  = "Men"
$ws.Range("C8").Value  = "Woman"

$ws.Range("B7").Value  = "Мужчины"
$ws.Range("B8").Value  = "Женщины"

$ws.Range("A7").Value  = "Эркектер"
$ws.Range("A8").Value  = "Аялдар"

$ws.Range("A22").Value = "Жаш курагы боюнча (айларда)"
$ws.Range("B22").Value = "По возрасту (в месяцах)"
